$d = $word.ActiveDocument
$t = $d.Tables(1)
$row = $t.Rows(3)

$startCell = $row.Cells(2)
$startCell.Range.Text = "00:00:10:00"

$stopCell = $row.Cells(3)
$stopCell.Range.Text = "00:00:20"
